$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column stays text (matches source workbook formatting)
$ws.Range("D2:D51").NumberFormat = "@"

$data = @(
    ,@('Bitcoin','https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc','28.970.57','  -1.91%  ')
    ,@('Ethereum','https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth','1.906.48','  -4.29%  ')
    ,@('TetherUSD','https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt','1.005','  +0.25%  ')
    ,@('BNB','https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb','325.10','  -0.02%  ')
    ,@('USDC','https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc','1.003','  +0.21%  ')
    ,@('XRP','https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp','0.4601','  -1.72%  ')
    ,@('Cardano','https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada','0.3823','  -3.15%  ')
    ,@('Dogecoin','https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge','0.07738','  -2.66%  ')
    ,@('Polygon','https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic','0.9829','  -1.86%  ')
    ,@('Solana','https://coinranking.com/coin/zNZHO_Sjf+solana-sol','22.09','  -3.77%  ')
    ,@('WrappedEther','https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth','1.940.74','  -3.04%  ')
    ,@('Chainlink','https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link','7.000','  -3.83%  ')
    ,@('Polkadot','https://coinranking.com/coin/25W7FG7om+polkadot-dot','5.688','  -3.09%  ')
    ,@('TRON','https://coinranking.com/coin/qUhEFk1I61atv+tron-trx','0.07042','  -1.47%  ')
    ,@('BinanceUSD','https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd','1.004','  +0.13%  ')
    ,@('Litecoin','https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc','84.26','  -5.04%  ')
    ,@('ShibaInu','https://coinranking.com/coin/xz24e0BjL+shibainu-shib','0.000009577','  -3.68%  ')
    ,@('Avalanche','https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax','16.76','  -3.58%  ')
    ,@('Dai','https://coinranking.com/coin/MoTuySvg7+dai-dai','1.004','  +0.20%  ')
    ,@('WrappedBTC','https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc','28.983.51','  -2.33%  ')
    ,@('Uniswap','https://coinranking.com/coin/_H5FVG9iW+uniswap-uni','5.342','  -3.41%  ')
    ,@('Cosmos','https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom','10.96','  -2.86%  ')
    ,@('WrappedliquidstakedEther2.0','https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth','2.170.79','  -3.14%  ')
    ,@('Toncoin','https://coinranking.com/coin/67YlI0K1b+toncoin-ton','2.077','  -1.29%  ')
    ,@('Monero','https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr','156.73','  -0.66%  ')
    ,@('EthereumClassic','https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc','19.17','  -2.52%  ')
    ,@('InternetComputer(DFINITY)','https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp','5.618','  -5.81%  ')
    ,@('BitcoinCash','https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch','117.91','  -2.06%  ')
    ,@('LidoDAOToken','https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo','1.831','  -6.58%  ')
    ,@('Stellar','https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm','0.09265','  -1.87%  ')
    ,@('ImmutableX','https://coinranking.com/coin/Z96jIvLU7+immutablex-imx','0.8621','  -4.63%  ')
    ,@('Filecoin','https://coinranking.com/coin/ymQub4fuB+filecoin-fil','5.115','  -2.67%  ')
    ,@('ARBITRUM','https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb','1.253','  -7.02%  ')
    ,@('HuobiToken','https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht','3.018','  -5.18%  ')
    ,@('Hedera','https://coinranking.com/coin/jad286TjB+hedera-hbar','0.05719','  -2.02%  ')
    ,@('TrustWalletToken','https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt','1.144','  -2.92%  ')
    ,@('Frax','https://coinranking.com/coin/KfWtaeV1W+frax-frax','1.003','  +0.13%  ')
    ,@('VeChain','https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet','0.02043','  -3.66%  ')
    ,@('FraxShare','https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs','7.514','  -4.55%  ')
    ,@('TheSandbox','https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand','0.5545','  -3.61%  ')
    ,@('Algorand','https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo','0.1759','  -3.86%  ')
    ,@('Aptos','https://coinranking.com/coin/HGYj5JCv5+aptos-apt','9.328','  -4.89%  ')
    ,@('MXToken','https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx','2.760','  +2.70%  ')
    ,@('Decentraland','https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana','0.5222','  -2.76%  ')
    ,@('EnergySwap','https://coinranking.com/coin/SbWqqTui-+energyswap-ens','11.29','  -6.67%  ')
    ,@('RenderToken','https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr','2.107','  -3.31%  ')
    ,@('PEPE','https://coinranking.com/coin/03WI8NQPF+pepe-pepe','0.000002634','  -19.58%  ')
    ,@('Cronos','https://coinranking.com/coin/65PHZTpmE55b+cronos-cro','0.06820','  -1.81%  ')
    ,@('Quant','https://coinranking.com/coin/bauj_21eYVwso+quant-qnt','112.16','  -1.84%  ')
    ,@('NEARProtocol','https://coinranking.com/coin/DCrsaMv68+nearprotocol-near','1.779','  -4.68%  ')
)

$arr = New-Object "object[,]" 50,4
for ($i = 0; $i -lt 50; $i++) {
    for ($j = 0; $j -lt 4; $j++) {
        $arr[$i, $j] = $data[$i][$j]
    }
}

$ws.Range("B2:E51").Value = $arr
